{"js": "const body = context.document.body;\n\nconst replacements = [\n  { find: \"Gris : #EBEBEB\", append: \"  un 60% de la p\u00e1gina, el fondo y lo menos important.\" },\n  { find: \"Naranja: #F5A25D\", append: \" un 10% para aquello que llama a la acci\u00f3n.\" },\n  { find: \"Rojo: #FA7F72\", append: \" solo para peque\u00f1os detalles o peque\u00f1os bloques de texto\" },\n  { find: \"Azul: #389393\", append: \" para el footer\" },\n];\n\nfor (const { find, append } of replacements) {\n  const results = body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"text\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${find}`);\n  }\n\n  const range = results.items[0];\n  range.insertText(find + append, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Find = \"Gris : #EBEBEB\"; Replace = \"Gris : #EBEBEB  un 60% de la p\u00e1gina, el fondo y lo menos important.\" },\n    @{ Find = \"Naranja: #F5A25D\"; Replace = \"Naranja: #F5A25D un 10% para aquello que llama a la acci\u00f3n.\" },\n    @{ Find = \"Rojo: #FA7F72\"; Replace = \"Rojo: #FA7F72 solo para peque\u00f1os detalles o peque\u00f1os bloques de texto\" },\n    @{ Find = \"Azul: #389393\"; Replace = \"Azul: #389393 para el footer\" }\n)\n\nforeach ($r in $replacements) {\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Text = $r.Find\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.Forward = $true\n    $find.Wrap = 0\n\n    if ($find.Execute()) {\n        $rng.Text = $r.Replace\n    }\n}\n"}
